# Generate Report for Handoff
# Moves the localization status from "In Translation" to "Ready for handoff"
# and refreshes the associated generate/handoff timestamps, on all three
# sheets (Overview, zh-cn, de-de). The "Status" columns widen slightly to
# fit the new, longer text, matching Excel's own recompute of those columns.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status: "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date (Overview) / Latest Handback DateTime (de-de):
# 2016-08-24 06:58:52 -> 2016-08-24 06:59:38
$wsOverview.Range("G2").Value = "2016-08-24 06:59:38"
$wsDeDe.Range("H2").Value = "2016-08-24 06:59:38"

# Latest Handoff Datetime (zh-cn): 2016-08-24 06:58:47 -> 2016-08-24 06:59:33
$wsZhCn.Range("H2").Value = "2016-08-24 06:59:33"

# Status column widths grow to fit "Ready for handoff" (was sized for "In Translation")
$wsOverview.Columns("E:F").ColumnWidth = 16.33
$wsZhCn.Columns("C:C").ColumnWidth = 16.33
$wsDeDe.Columns("C:C").ColumnWidth = 16.33
